$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns and apply municipality/state text normalization
# (capitalize connector words De/Del/La/Los/El/Y, plus a couple of special-case fixes)
$changes = @{
    'A1' = 'mx_state'
    'B1' = 'mx_municipality'
    'C1' = 'n_matriculas'
    'D1' = 'pct_matriculas'
    'B8' = 'Pabellón De Arteaga'
    'B9' = 'Rincón De Romos'
    'B10' = 'San Francisco De Los Romo'
    'B11' = 'San José De Gracia'
    'B41' = 'Mazapa De Madero'
    'B44' = 'Ocozocoautla De Espinosa'
    'B71' = 'Coyame Del Sotol'
    'B80' = 'Guadalupe Y Calvo'
    'B82' = 'Hidalgo Del Parral'
    'B103' = 'Valle De Zaragoza'
    'B118' = 'San Juan De Sabinas'
    'A130' = 'Ciudad De México'
    'B134' = 'Cuajimalpa De Morelos'
    'B147' = 'Coneto De Comonfort'
    'B159' = 'Nombre De Dios'
    'B163' = 'Pánuco De Coronado'
    'B169' = 'San Juan De Guadalupe'
    'B170' = 'San Juan Del Río'
    'B171' = 'San Luis Del Cordero'
    'B172' = 'San Pedro Del Gallo'
    'A181' = 'Estado De México'
    'B181' = 'Acambay De Ruíz Castañeda'
    'B186' = 'Atizapán De Zaragoza'
    'B194' = 'Ecatepec De Morelos'
    'B201' = 'Naucalpan De Juárez'
    'B206' = 'San Felipe Del Progreso'
    'B211' = 'Tenango Del Valle'
    'B214' = 'Tlalnepantla De Baz'
    'B217' = 'Valle De Bravo'
    'B218' = 'Villa De Allende'
    'B219' = 'Villa Del Carbón'
    'A225' = 'Guanajuato'
    'B228' = 'Apaseo El Grande'
    'B234' = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    'B238' = 'Jaral Del Progreso'
    'B248' = 'San Diego De La Unión'
    'B250' = 'San Francisco Del Rincón'
    'B252' = 'San Luis De La Paz'
    'B253' = 'Santa Cruz De Juventino Rosas'
    'B257' = 'Valle De Santiago'
    'B263' = 'Acapulco De Juárez'
    'B264' = 'Ajuchitlán Del Progreso'
    'B267' = 'Atoyac De Álvarez'
    'B268' = 'Ayutla De Los Libres'
    'B271' = 'Chilapa De Álvarez'
    'B272' = 'Chilpancingo De Los Bravo'
    'B274' = 'Coyuca De Catalán'
    'B277' = 'Cutzamala De Pinzón'
    'B282' = 'Iguala De La Independencia'
    'B283' = 'Zihuatanejo De Azueta'
    'B293' = 'Taxco De Alarcón'
    'B294' = 'Técpan De Galeana'
    'B296' = 'Tepecoacuilco De Trujano'
    'B298' = 'Tlapa De Comonfort'
    'B313' = 'Jacala De Ledezma'
    'B317' = 'Mixquiahuala De Juárez'
    'B318' = 'Molango De Escamilla'
    'B320' = 'Nopala De Villagrán'
    'B321' = 'Pachuca De Soto'
    'B322' = 'Progreso De Obregón'
    'B324' = 'Santiago De Anaya'
    'B328' = 'Tepeji Del Río De Ocampo'
    'B329' = 'Tezontepec De Aldama'
    'B332' = 'Tula De Allende'
    'B333' = 'Tulancingo De Bravo'
    'B334' = 'Zacualtipán De Ángeles'
    'B339' = 'Atotonilco El Alto'
    'B340' = 'Autlán De Navarro'
    'B352' = 'Encarnación De Díaz'
    'B357' = 'Ixtlahuacán De Los Membrillos'
    'B360' = 'Jilotlán De Los Dolores'
    'B362' = 'Lagos De Moreno'
    'B367' = 'Ojuelos De Jalisco'
    'B371' = 'San Juan De Los Lagos'
    'B373' = 'San Martín De Bolaños'
    'B375' = 'Santa María De Los Ángeles'
    'B378' = 'Tamazula De Gordiano'
    'B382' = 'Teocuitatlán De Corona'
    'B383' = 'Tepatitlán De Morelos'
    'B388' = 'Unión De San Antonio'
    'B390' = 'Yahualica De González Gallo'
    'B392' = 'Zapotitlán De Vadillo'
    'B393' = 'Zapotlán El Grande'
    'B406' = 'Coalcomán De Vázquez Pallares'
    'B442' = 'Tiquicheo De Nicolás Romero'
    'B457' = 'Coatlán Del Río'
    'B467' = 'Tetela Del Volcán'
    'B472' = 'Bahía De Banderas'
    'B492' = 'Mier Y Noriega'
    'B493' = 'Montemorelos'
    'B495' = 'San Nicolás De Los Garza'
    'B499' = 'Acatlán De Pérez Figueroa'
    'B501' = 'Ayoquezco De Aldama'
    'B502' = 'Ciénega De Zimatlán'
    'B504' = 'El Barrio De La Soledad'
    'B505' = 'Heroica Ciudad De Ejutla De Crespo'
    'B506' = 'Heroica Ciudad De Huajuapan De León'
    'B507' = 'Heroica Ciudad De Juchitán De Zaragoza'
    'B509' = 'Mariscala De Juárez'
    'B510' = 'Mártires De Tacubaya'
    'B512' = 'Miahuatlán De Porfirio Díaz'
    'B514' = 'Nejapa De Madero'
    'B515' = 'Oaxaca De Juárez'
    'B516' = 'Ocotlán De Morelos'
    'B517' = 'Putla Villa De Guerero'
    'B519' = 'San Agustín De Las Juntas'
    'B523' = 'San Dionisio Del Mar'
    'B558' = 'Santo Domingo De Morelos'
    'B563' = 'Tepelmeme Villa De Morelos'
    'B564' = 'Villa De Etla'
    'B565' = 'Villa De Tututepec De Melchor Ocampo'
    'B566' = 'Zimatlán De Álvarez'
    'B576' = 'Cuetzalan Del Progreso'
    'B579' = 'Huehuetlán El Chico'
    'B583' = 'Izúcar De Matamoros'
    'B586' = 'Palmar De Bravo'
    'B592' = 'San Salvador El Seco'
    'B593' = 'Tecali De Herrera'
    'B596' = 'Tepexi De Rodríguez'
    'B602' = 'Totoltepec De Guerero'
    'B603' = 'Tuzamapan De Galeana'
    'B609' = 'Amealco De Bonfil'
    'B611' = 'Cadereyta De Montes'
    'B615' = 'Jalpan De Serra'
    'B616' = 'Landa De Matamoros'
    'B617' = 'Pinal De Amoles'
    'B619' = 'San Juan Del Río'
    'B631' = 'Ciudad Del Maíz'
    'B641' = 'Mexquitic De Carmona'
    'B651' = 'Santa María Del Río'
    'B652' = 'Soledad De Graciano Sánchez'
    'B660' = 'Villa De Arriaga'
    'B661' = 'Villa De Guadalupe'
    'B662' = 'Villa De Ramos'
    'B663' = 'Villa De Reyes'
    'B687' = 'Nacozari De García'
    'B699' = 'Jalpa De Méndez'
    'B722' = 'Soto La Marina'
    'B730' = 'Ixtacuixtla De Mariano Matamoros'
    'B731' = 'Nanacamilpa De Mariano Arista'
    'B742' = 'Castillo De Teayo'
    'B750' = 'Cosamaloapan De Carpio'
    'B751' = 'Cosautlán De Carvajal'
    'B758' = 'Hueyapan De Ocampo'
    'B759' = 'Ignacio De La Llave'
    'B762' = 'Ixhuacán De Los Reyes'
    'B763' = 'Ixhuatlán De Madero'
    'B770' = 'Martínez De La Torre'
    'B780' = 'Poza Rica De Hidalgo'
    'B786' = 'Sayula De Alemán'
    'B790' = 'Tatahuicapan De Juárez'
    'B816' = 'Concepción Del Oro'
    'B834' = 'Nochistlán De Mejía'
    'B835' = 'Noria De Ángeles'
    'B844' = 'Tlaltenango De Sánchez Román'
    'B847' = 'Villa De Cos'
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}

# Remove the trailing metadata/footnote rows (854-858) that are no longer part of the table
$ws.Range("A854:D858").EntireRow.Delete()

